$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.082.11"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.789.57"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'227.03"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'32.26"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  +3.82%  "
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").Value = "2.046.65"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "'11.34"
$ws.Range("E13").Value = "  +5.32%  "
$ws.Range("D14").Value = "1.790.12"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E15").Value = "  +0.29%  "
$ws.Range("D16").Value = "34.079.25"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").Value = "'243.82"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").Value = "'162.02"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").Value = "'7.20"
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'3.61"
$ws.Range("E33").Value = "  +3.61%  "
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "1.411.20"
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("E38").Value = "  +7.99%  "
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'80.47"
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.921"
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("D44").Value = "'13.34"
$ws.Range("E44").Value = "  +8.83%  "
$ws.Range("D45").Value = "0.0₆0139"
$ws.Range("E45").Value = "  -5.60%  "
$ws.Range("D46").Value = "'0.0507"
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  +3.36%  "
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "1.947.33"
$ws.Range("E50").Value = "  +0.51%  "
$ws.Range("E51").Value = "  +0.04%  "
